$wb = $excel.ActiveWorkbook

# --- Sheet "REPCA1": insert a new column "Kc" after "Xc" (before "emax"), i.e. before column T ---
$repca = $wb.Worksheets.Item("REPCA1")
$repca.Columns("T").Insert()

# Header for new column T = "Kc" (match the bold/bordered header style used
# by the rest of row 1, e.g. S1)
$repca.Range("T1").Value = "Kc"
$repca.Range("T1").Font.Bold = $true
$repca.Range("T1").HorizontalAlignment = -4108
$repca.Range("T1").VerticalAlignment = -4160
$repca.Range("T1").Borders.LineStyle = 1

# Data value for new column T, row 2 = 1 (Kc)
$repca.Range("T2").Value = 1

# --- Other data updates on REPCA1 row 2 ---
$repca.Range("I2").Value = 1        # VCFlag: 0 -> 1
$repca.Range("J2").Value = 1        # RefFlag: 0 -> 1
$repca.Range("Q2").Value = 0.98     # Vfrz: 0.8 -> 0.98

# dbd1/dbd2 now sit at W2/X2 after the column insert (previously V2/W2)
$repca.Range("W2").Value = -0.02    # dbd1: -0.1 -> -0.02
$repca.Range("X2").Value = 0.02     # dbd2: 0.1 -> 0.02

# Pmin now sits at AI2 after the insert (previously AH2)
$repca.Range("AI2").Value = -999    # Pmin: 0 -> -999

# Ddn/Dup now sit at AK2/AL2 after the insert (previously AJ2/AK2)
$repca.Range("AK2").Value = 10      # Ddn: 0.05 -> 10
$repca.Range("AL2").Value = 10      # Dup: 0.05 -> 10

# Update the view: scroll right (frozen header row stays, columns scroll so
# column S is the first visible one) and select AL3 in the scrollable pane
$repca.Activate()
$repca.Range("AL3").Select()
$excel.ActiveWindow.ScrollColumn = 19

# --- Sheet "Toggler": enable the second toggler row (u: 0 -> 1) ---
$toggler = $wb.Worksheets.Item("Toggler")
$toggler.Range("C3").Value = 1
$toggler.Activate()
$toggler.Range("G4").Select()

# REPCA1 should be the final active/selected sheet (tabSelected) per the workbook view changes
$repca.Activate()
